$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.017265498297682
$ws.Cells.Item(2, 4).Value = 1.022552066106101
$ws.Cells.Item(2, 5).Value = 1.018687445170994
$ws.Cells.Item(2, 9).Value = 1.026631938829507
$ws.Cells.Item(2, 10).Value = 1.022480596013615
$ws.Cells.Item(2, 11).Value = 1.025386111524837
$ws.Cells.Item(2, 12).Value = 1.021532920243654
$ws.Cells.Item(2, 14).Value = 1.023932635740118
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.01812685219142
$ws.Cells.Item(3, 4).Value = 1.023145628004424
$ws.Cells.Item(3, 5).Value = 1.019415355757227
$ws.Cells.Item(3, 9).Value = 1.026713740711176
$ws.Cells.Item(3, 10).Value = 1.02297822480766
$ws.Cells.Item(3, 11).Value = 1.025786946181015
$ws.Cells.Item(3, 12).Value = 1.022066909360582
$ws.Cells.Item(3, 14).Value = 1.024430971224129
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.018684726435612
$ws.Cells.Item(4, 4).Value = 1.023529876056189
$ws.Cells.Item(4, 5).Value = 1.019887214873861
$ws.Cells.Item(4, 9).Value = 1.02676534204749
$ws.Cells.Item(4, 10).Value = 1.023300126374124
$ws.Cells.Item(4, 11).Value = 1.026045790930264
$ws.Cells.Item(4, 12).Value = 1.022412621914033
$ws.Cells.Item(4, 14).Value = 1.02475332992774
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.018919379882687
$ws.Cells.Item(5, 4).Value = 1.023691453723528
$ws.Cells.Item(5, 5).Value = 1.020085786700008
$ws.Cells.Item(5, 9).Value = 1.026786716192964
$ws.Cells.Item(5, 10).Value = 1.023435429136954
$ws.Cells.Item(5, 11).Value = 1.026154482951261
$ws.Cells.Item(5, 12).Value = 1.022558002591647
$ws.Cells.Item(5, 14).Value = 1.024888824836012
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.018958786403865
$ws.Cells.Item(6, 4).Value = 1.023718585604795
$ws.Cells.Item(6, 5).Value = 1.020119139591606
$ws.Cells.Item(6, 9).Value = 1.026790286277134
$ws.Cells.Item(6, 10).Value = 1.023458145582411
$ws.Cells.Item(6, 11).Value = 1.0261727253743
$ws.Cells.Item(6, 12).Value = 1.022582415135501
$ws.Cells.Item(6, 14).Value = 1.024911573541429
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.018687861404622
$ws.Cells.Item(7, 4).Value = 1.023532034910858
$ws.Cells.Item(7, 5).Value = 1.019889867407854
$ws.Cells.Item(7, 9).Value = 1.026765628904371
$ws.Cells.Item(7, 10).Value = 1.023301934394784
$ws.Cells.Item(7, 11).Value = 1.026047243776356
$ws.Cells.Item(7, 12).Value = 1.022414564330069
$ws.Cells.Item(7, 14).Value = 1.024755140515996
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.017556488144713
$ws.Cells.Item(8, 4).Value = 1.02275262604363
$ws.Cells.Item(8, 5).Value = 1.018933268414881
$ws.Cells.Item(8, 9).Value = 1.026659859134844
$ws.Cells.Item(8, 10).Value = 1.022648791406441
$ws.Cells.Item(8, 11).Value = 1.025521682776714
$ws.Cells.Item(8, 12).Value = 1.021713345046232
$ws.Cells.Item(8, 14).Value = 1.024101069989694
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.015566919621536
$ws.Cells.Item(9, 4).Value = 1.02138061976885
$ws.Cells.Item(9, 5).Value = 1.017254225255249
$ws.Cells.Item(9, 9).Value = 1.026463322752349
$ws.Cells.Item(9, 10).Value = 1.021497180595354
$ws.Cells.Item(9, 11).Value = 1.024591637169025
$ws.Cells.Item(9, 12).Value = 1.020479191775503
$ws.Cells.Item(9, 14).Value = 1.022947823759164
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.01424335853714
$ws.Cells.Item(10, 4).Value = 1.020467006681796
$ws.Cells.Item(10, 5).Value = 1.016139411385556
$ws.Cells.Item(10, 9).Value = 1.026325509281334
$ws.Cells.Item(10, 10).Value = 1.020729054047554
$ws.Cells.Item(10, 11).Value = 1.023969036539913
$ws.Cells.Item(10, 12).Value = 1.019657503324466
$ws.Cells.Item(10, 14).Value = 1.02217860638356
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.013670929399835
$ws.Cells.Item(11, 4).Value = 1.020071675936195
$ws.Cells.Item(11, 5).Value = 1.015657785902545
$ws.Cells.Item(11, 9).Value = 1.026264232509397
$ws.Cells.Item(11, 10).Value = 1.020396371513466
$ws.Cells.Item(11, 11).Value = 1.023698850560095
$ws.Cells.Item(11, 12).Value = 1.019301976710557
$ws.Cells.Item(11, 14).Value = 1.021845451402114
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.013458407640179
$ws.Cells.Item(12, 4).Value = 1.019924874874043
$ws.Cells.Item(12, 5).Value = 1.015479055275208
$ws.Cells.Item(12, 9).Value = 1.026241231460614
$ws.Cells.Item(12, 10).Value = 1.020272787892932
$ws.Cells.Item(12, 11).Value = 1.023598403272573
$ws.Cells.Item(12, 12).Value = 1.019169960238604
$ws.Cells.Item(12, 14).Value = 1.021721692278664
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.013503989552584
$ws.Cells.Item(13, 4).Value = 1.019956362252785
$ws.Cells.Item(13, 5).Value = 1.015517386036885
$ws.Cells.Item(13, 9).Value = 1.026246176119496
$ws.Cells.Item(13, 10).Value = 1.020299297453929
$ws.Cells.Item(13, 11).Value = 1.023619953525972
$ws.Cells.Item(13, 12).Value = 1.019198276299677
$ws.Cells.Item(13, 14).Value = 1.021748239486279
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.013653360156396
$ws.Cells.Item(14, 4).Value = 1.020059540439393
$ws.Cells.Item(14, 5).Value = 1.015643008565644
$ws.Cells.Item(14, 9).Value = 1.026262336130297
$ws.Cells.Item(14, 10).Value = 1.020386156259795
$ws.Cells.Item(14, 11).Value = 1.023690549347243
$ws.Cells.Item(14, 12).Value = 1.01929106332415
$ws.Cells.Item(14, 14).Value = 1.021835221641612
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.013745406129608
$ws.Cells.Item(15, 4).Value = 1.020123117596395
$ws.Cells.Item(15, 5).Value = 1.01572043087171
$ws.Cells.Item(15, 9).Value = 1.026272261046527
$ws.Cells.Item(15, 10).Value = 1.020439671482496
$ws.Cells.Item(15, 11).Value = 1.023734034117037
$ws.Cells.Item(15, 12).Value = 1.019348238071642
$ws.Cells.Item(15, 14).Value = 1.021888812862067
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.014281363262366
$ws.Cells.Item(16, 4).Value = 1.020493249322909
$ws.Cells.Item(16, 5).Value = 1.01617139853525
$ws.Cells.Item(16, 9).Value = 1.02632954230877
$ws.Cells.Item(16, 10).Value = 1.02075113154558
$ws.Cells.Item(16, 11).Value = 1.023986955454125
$ws.Cells.Item(16, 12).Value = 1.019681104260414
$ws.Cells.Item(16, 14).Value = 1.022200715234166
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.014617738549771
$ws.Cells.Item(17, 4).Value = 1.020725496704847
$ws.Cells.Item(17, 5).Value = 1.016454573425467
$ws.Cells.Item(17, 9).Value = 1.026365044720435
$ws.Cells.Item(17, 10).Value = 1.02094648207367
$ws.Cells.Item(17, 11).Value = 1.024145447752364
$ws.Cells.Item(17, 12).Value = 1.019889975671418
$ws.Cells.Item(17, 14).Value = 1.022396343182413
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.014814006221524
$ws.Cells.Item(18, 4).Value = 1.020860988576441
$ws.Cells.Item(18, 5).Value = 1.016619850201283
$ws.Cells.Item(18, 9).Value = 1.026385598002435
$ws.Cells.Item(18, 10).Value = 1.021060419053698
$ws.Cells.Item(18, 11).Value = 1.024237836053597
$ws.Cells.Item(18, 12).Value = 1.02001183278058
$ws.Cells.Item(18, 14).Value = 1.022510441966021
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.014880939516209
$ws.Cells.Item(19, 4).Value = 1.020907192140422
$ws.Cells.Item(19, 5).Value = 1.016676223175987
$ws.Cells.Item(19, 9).Value = 1.02639257988528
$ws.Cells.Item(19, 10).Value = 1.021099267277473
$ws.Cells.Item(19, 11).Value = 1.024269328270678
$ws.Cells.Item(19, 12).Value = 1.020053387264761
$ws.Cells.Item(19, 14).Value = 1.02254934535873
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.014581641862634
$ws.Cells.Item(20, 4).Value = 1.020700576066087
$ws.Cells.Item(20, 5).Value = 1.016424180498347
$ws.Cells.Item(20, 9).Value = 1.02636125164163
$ws.Cells.Item(20, 10).Value = 1.020925523595526
$ws.Cells.Item(20, 11).Value = 1.024128448974368
$ws.Cells.Item(20, 12).Value = 1.019867563051685
$ws.Cells.Item(20, 14).Value = 1.022375354940826
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.013609371355954
$ws.Cells.Item(21, 4).Value = 1.020029155850373
$ws.Cells.Item(21, 5).Value = 1.015606011250593
$ws.Cells.Item(21, 9).Value = 1.026257584033367
$ws.Cells.Item(21, 10).Value = 1.020360578786153
$ws.Cells.Item(21, 11).Value = 1.02366976306111
$ws.Cells.Item(21, 12).Value = 1.019263738689119
$ws.Cells.Item(21, 14).Value = 1.021809607845023
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.012998668238277
$ws.Cells.Item(22, 4).Value = 1.01960725341642
$ws.Cells.Item(22, 5).Value = 1.01509256011088
$ws.Cells.Item(22, 9).Value = 1.026191015112243
$ws.Cells.Item(22, 10).Value = 1.020005315372057
$ws.Cells.Item(22, 11).Value = 1.023380859468836
$ws.Cells.Item(22, 12).Value = 1.018884333699299
$ws.Cells.Item(22, 14).Value = 1.02145383991613
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.013322355952491
$ws.Cells.Item(23, 4).Value = 1.019830887934359
$ws.Cells.Item(23, 5).Value = 1.015364658260375
$ws.Cells.Item(23, 9).Value = 1.02622643599935
$ws.Cells.Item(23, 10).Value = 1.020193652499694
$ws.Cells.Item(23, 11).Value = 1.023534060586481
$ws.Cells.Item(23, 12).Value = 1.019085439931744
$ws.Cells.Item(23, 14).Value = 1.021642444504091
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.014597952199079
$ws.Cells.Item(24, 4).Value = 1.020711836550241
$ws.Cells.Item(24, 5).Value = 1.016437913427641
$ws.Cells.Item(24, 9).Value = 1.026362966048664
$ws.Cells.Item(24, 10).Value = 1.020934993853528
$ws.Cells.Item(24, 11).Value = 1.024136130168495
$ws.Cells.Item(24, 12).Value = 1.019877690269792
$ws.Cells.Item(24, 14).Value = 1.022384838647681
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.016080780773281
$ws.Cells.Item(25, 4).Value = 1.021735137773461
$ws.Cells.Item(25, 5).Value = 1.017687504281337
$ws.Cells.Item(25, 9).Value = 1.026515331652794
$ws.Cells.Item(25, 10).Value = 1.021794973300989
$ws.Cells.Item(25, 11).Value = 1.024832535031713
$ws.Cells.Item(25, 12).Value = 1.020798065560595
$ws.Cells.Item(25, 14).Value = 1.023246039364598
